$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New column I: "End date" header + "14-MAR-2014" for the 3 data rows.
# Style is copied (format-only paste) from column H so the new cells
# reuse H's existing style indices instead of minting new ones, and the
# "14-MAR-2014" text must land as literal text (not an auto-parsed date
# serial) -- so number format is forced to Text before the value is
# typed, then the H-column format is re-pasted on top to restore the
# original (General) number format bit without touching the already
# text-typed cell value.
# ---------------------------------------------------------------------

$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "End date"

$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "14-MAR-2014"
$ws.Range("H2").Copy()
$ws.Range("I2").PasteSpecial(-4122)

$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "14-MAR-2014"
$ws.Range("H3").Copy()
$ws.Range("I3").PasteSpecial(-4122)

$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "14-MAR-2014"
$ws.Range("H4").Copy()
$ws.Range("I4").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# New col I is 12 "screen" units wide; the ColumnWidth property is in
# character-width units, offset from the stored XML width by the default
# font's padding (~0.8333), so back that offset out to land on width=12.
$ws.Columns.Item(9).ColumnWidth = 11.166666666666666

# Row heights for data rows (2-4) increase to accommodate wrapped text
$ws.Rows.Item(2).RowHeight = 26.25
$ws.Rows.Item(3).RowHeight = 26.25
$ws.Rows.Item(4).RowHeight = 26.25

# Update view: scroll so column C is left-most visible, select J7
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("J7").Select()
